$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.365.97'
$ws.Range('E2').Value = '  +12.74%  '
$ws.Range('D3').Value = '1.824.46'
$ws.Range('E3').Value = '  +9.04%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'229.64"
$ws.Range('E5').Value = '  +4.63%  '
$ws.Range('D6').Value = "'0.548"
$ws.Range('E6').Value = '  +3.34%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = "'31.62"
$ws.Range('E8').Value = '  +6.65%  '
$ws.Range('D9').Value = "'47.16"
$ws.Range('E9').Value = '  +6.86%  '
$ws.Range('E10').Value = '  +7.34%  '
$ws.Range('D11').Value = "'0.0673"
$ws.Range('E11').Value = '  +5.48%  '
$ws.Range('D12').Value = "'0.0930"
$ws.Range('E12').Value = '  +2.76%  '
$ws.Range('E13').Value = '  +9.04%  '
$ws.Range('D14').Value = '1.830.43'
$ws.Range('E14').Value = '  +9.53%  '
$ws.Range('E15').Value = '  +5.55%  '
$ws.Range('E16').Value = '  +2.20%  '
$ws.Range('D17').Value = '34.279.95'
$ws.Range('E17').Value = '  +12.32%  '
$ws.Range('E18').Value = '  +7.72%  '
$ws.Range('D19').Value = "'69.77"
$ws.Range('E19').Value = '  +5.13%  '
$ws.Range('D20').Value = "'258.74"
$ws.Range('E20').Value = '  +6.67%  '
$ws.Range('E21').Value = '  +4.77%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').Value = "'10.62"
$ws.Range('E23').Value = '  +6.47%  '
$ws.Range('D24').Value = "'4.35"
$ws.Range('E24').Value = '  +1.83%  '
$ws.Range('D25').Value = "'2.22"
$ws.Range('E25').Value = '  +2.66%  '
$ws.Range('D26').Value = "'159.75"
$ws.Range('E26').Value = '  +0.88%  '
$ws.Range('E27').Value = '  +5.00%  '
$ws.Range('E28').Value = '  +7.41%  '
$ws.Range('E29').Value = '  +2.49%  '
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').Value = "'3.92"
$ws.Range('E31').Value = '  +12.97%  '
$ws.Range('E32').Value = '  +5.16%  '
$ws.Range('E33').Value = '  +6.32%  '
$ws.Range('E34').Value = '  +8.67%  '
$ws.Range('D35').Value = '1.550.20'
$ws.Range('E35').Value = '  +3.32%  '
$ws.Range('D36').Value = "'1.79"
$ws.Range('E36').Value = '  +1.82%  '
$ws.Range('E37').Value = '  +6.56%  '
$ws.Range('B38').Value = 'MinaProtocolToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina'
$ws.Range('D38').Value = "'1.30"
$ws.Range('E38').Value = '  +215.87%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').Value = "'0.635"
$ws.Range('E39').Value = '  +6.27%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.0190"
$ws.Range('E40').Value = '  +6.93%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = "'84.74"
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = "'2.82"
$ws.Range('E42').Value = '  +5.49%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = "'0.920"
$ws.Range('E43').Value = '  +9.74%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').Value = "'2.33"
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').Value = "'2.17"
$ws.Range('E45').Value = '  +10.06%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').Value = "'0.0527"
$ws.Range('E46').Value = '  +6.06%  '
$ws.Range('E47').Value = '  +4.76%  '
$ws.Range('E48').Value = '  +10.18%  '
$ws.Range('D49').Value = "'12.34"
$ws.Range('E49').Value = '  +27.47%  '
$ws.Range('E50').Value = '  +4.31%  '
$ws.Range('D51').Value = "'0.999"
$ws.Range('E51').Value = '  -0.12%  '
